# The document's header/footer logo pictures (the Pearson Edexcel logo in
# the footers and the BTec logo in the headers) were exported with
# duplicate/mismatched drawing names ("image2.png" used twice for the
# Pearson logo, "image1.jpg" used twice for the BTec logo). This renames
# each inline picture so the two logo types no longer collide:
#   - Pearson logo pictures (footers):  image2.png -> image1.png
#   - BTec logo pictures   (headers):   image1.jpg -> image2.jpg

$d = $word.ActiveDocument

for ($si = 1; $si -le $d.Sections.Count; $si++) {
    $sec = $d.Sections.Item($si)

    # Headers: rename the BTec logo picture(s) image1.jpg -> image2.jpg
    for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        if ($hdr.Exists) {
            $shapes = $hdr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }

    # Footers: rename the Pearson Edexcel logo picture(s) image2.png -> image1.png
    for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
        $ftr = $sec.Footers.Item($fi)
        if ($ftr.Exists) {
            $shapes = $ftr.Range.InlineShapes
            for ($i = 1; $i -le $shapes.Count; $i++) {
                $shp = $shapes.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}
